$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "shankar" row (row 2): email and username
$ws.Range("C2").Value = "shankar1222"
$ws.Range("B2").Value = "shankar72@gmail.com"

# Update the "prem" row (row 3): email and username
$ws.Range("C3").Value = "prem1222"
$ws.Range("B3").Value = "prem23@gmail.com"

# Update the selected cell to B2
$ws.Range("B2").Select()
